$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C3").Value = 8
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 4
$ws.Range("C9").Value = 5
$ws.Range("C10").Value = 4
$ws.Range("C11").Value = 2
$ws.Range("C12").Value = 3

# Row 13: update label and value
$ws.Range("B13").Value = "<delpa>"
$ws.Range("C13").Value = 10

# Row 14: update label and value
$ws.Range("B14").Value = "<thet>"
$ws.Range("C14").Value = 5

$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 9
